$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, formatted the same way as the other header cells (copy format from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column values
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
